# Process map template update: remove the "McKinsey" brand name from the
# footer / logo text boxes, e.g.
#   "McKinsey & Company"                                              -> " & Company"
#   "...specific permission of McKinsey & Company is strictly ..."    -> "...specific permission of & Company is strictly ..."
#
# These text boxes aren't on the (single) content slide itself - they live on
# the slide layouts ("masters" in the simplified PowerPoint UI) that are
# shared across the deck. Presentation.SlideMaster.CustomLayouts exposes the
# full, flattened set of layouts across every slide master/design in the
# file, so walk that collection and patch any run containing "McKinsey".

$p = $ppt.ActivePresentation

$oldDisclaimer = "Any use of this material without specific permission of McKinsey & Company is strictly prohibited"
$newDisclaimer = "Any use of this material without specific permission of & Company is strictly prohibited"

$oldLogo = "McKinsey & Company"
$newLogo = " & Company"

$layouts = $p.SlideMaster.CustomLayouts

for ($i = 1; $i -le $layouts.Count; $i++) {
    $layout = $layouts.Item($i)

    for ($j = 1; $j -le $layout.Shapes.Count; $j++) {
        $shp = $layout.Shapes.Item($j)

        if ($shp.HasTextFrame) {
            $tf = $shp.TextFrame
            if ($tf.HasText) {
                $tr = $tf.TextRange
                $text = $tr.Text

                if ($text.Contains($oldDisclaimer)) {
                    $tr.Text = $text.Replace($oldDisclaimer, $newDisclaimer)
                } elseif ($text.Contains($oldLogo)) {
                    $tr.Text = $text.Replace($oldLogo, $newLogo)
                } elseif ($text -like "*McKinsey*") {
                    # Fallback for any other occurrence: drop the bare word.
                    $tr.Text = $text.Replace("McKinsey", "")
                }
            }
        }
    }
}
